# Auto-applied updates to Exodus_Profits market-data sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 383155.38
$ws.Range("I137").Value = 1463.2646
$ws.Range("J137").Value = 3627538.2
$ws.Range("K137").Value = 4389.793799999999
$ws.Range("L137").Value = 10882614.6
$ws.Range("M137").Value = -1839.793799999999
$ws.Range("N137").Value = -10887714.6


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7581.3647
$ws.Range("I32").Value = 3742.4312
$ws.Range("K32").Value = 3742.4312
$ws.Range("M32").Value = -3455.4312

$ws.Range("H61").Value = 86040.836
$ws.Range("I61").Value = 1728.8
$ws.Range("J61").Value = 146263.72
$ws.Range("K61").Value = 1728.8
$ws.Range("L61").Value = 146263.72
$ws.Range("M61").Value = -1516.8
$ws.Range("N61").Value = -146687.72

$ws.Range("H74").Value = 26897.6
$ws.Range("I74").Value = 37084.035
$ws.Range("J74").Value = 3129.25
$ws.Range("K74").Value = 37084.035
$ws.Range("L74").Value = 3129.25
$ws.Range("M74").Value = -36210.035
$ws.Range("N74").Value = -4877.25

$ws.Range("H76").Value = 134997.4
$ws.Range("J76").Value = 134997.4
$ws.Range("L76").Value = 134997.4
$ws.Range("N76").Value = -135673.4

$ws.Range("H77").Value = 26897.6
$ws.Range("I77").Value = 37084.035
$ws.Range("J77").Value = 3129.25
$ws.Range("K77").Value = 185420.175
$ws.Range("L77").Value = 15646.25
$ws.Range("M77").Value = -181052.175
$ws.Range("N77").Value = -24382.25

$ws.Range("H79").Value = 134997.4
$ws.Range("J79").Value = 134997.4
$ws.Range("L79").Value = 134997.4
$ws.Range("N79").Value = -137337.4

$ws.Range("H132").Value = 1704.7347
$ws.Range("I132").Value = 1672.5869
$ws.Range("K132").Value = 5017.7607
$ws.Range("M132").Value = -2487.7607

$ws.Range("H136").Value = 86040.836
$ws.Range("I136").Value = 1728.8
$ws.Range("J136").Value = 146263.72
$ws.Range("K136").Value = 5186.4
$ws.Range("L136").Value = 438791.16
$ws.Range("M136").Value = -2636.4
$ws.Range("N136").Value = -443891.16

$ws.Range("H138").Value = 112595.336
$ws.Range("J138").Value = 112595.336
$ws.Range("L138").Value = 112595.336
$ws.Range("N138").Value = -122875.336


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H63").Value = 9999
$ws.Range("I63").Value = 9999
$ws.Range("K63").Value = 9999
$ws.Range("M63").Value = -9313

$ws.Range("H66").Value = 9999
$ws.Range("I66").Value = 9999
$ws.Range("K66").Value = 29997
$ws.Range("M66").Value = -26565

$ws.Range("H105").Value = 50499.24
$ws.Range("I105").Value = 144997
$ws.Range("K105").Value = 144997
$ws.Range("M105").Value = -143250

$ws.Range("H134").Value = 2046.3636
$ws.Range("I134").Value = 1698.4736
$ws.Range("J134").Value = 4249.6665
$ws.Range("K134").Value = 5095.4208
$ws.Range("L134").Value = 12748.9995
$ws.Range("M134").Value = -2560.4208
$ws.Range("N134").Value = -17818.9995

$ws.Range("H140").Value = 65915.875
$ws.Range("J140").Value = 65915.875
$ws.Range("L140").Value = 65915.875
$ws.Range("N140").Value = -76275.875

$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6392.625
$ws.Range("I7").Value = 7766.385
$ws.Range("J7").Value = 5452.684
$ws.Range("K7").Value = 7766.385
$ws.Range("L7").Value = 5452.684
$ws.Range("M7").Value = -7653.385
$ws.Range("N7").Value = -5678.684

$ws.Range("H31").Value = 5410.1665
$ws.Range("I31").Value = 1886.8572
$ws.Range("K31").Value = 1886.8572
$ws.Range("M31").Value = -1591.8572

$ws.Range("H34").Value = 5410.1665
$ws.Range("I34").Value = 1886.8572
$ws.Range("K34").Value = 1886.8572
$ws.Range("M34").Value = -1684.8572

$ws.Range("H58").Value = 2467.7273
$ws.Range("J58").Value = 4399.3335
$ws.Range("L58").Value = 4399.3335
$ws.Range("N58").Value = -4805.3335

$ws.Range("H132").Value = 1096.5
$ws.Range("I132").Value = 558.125
$ws.Range("K132").Value = 1674.375
$ws.Range("M132").Value = 855.625

$ws.Range("H136").Value = 2467.7273
$ws.Range("J136").Value = 4399.3335
$ws.Range("L136").Value = 13198.0005
$ws.Range("N136").Value = -18298.0005


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 559.0833
$ws.Range("I17").Value = 11
$ws.Range("K17").Value = 33
$ws.Range("M17").Value = 136

$ws.Range("H86").Value = 3636.5264
$ws.Range("J86").Value = 3788.611
$ws.Range("L86").Value = 11365.833
$ws.Range("N86").Value = -13737.833

$ws.Range("H89").Value = 3636.5264
$ws.Range("J89").Value = 3788.611
$ws.Range("L89").Value = 34097.499
$ws.Range("N89").Value = -45953.499

$ws.Range("H97").Value = 213.57143
$ws.Range("J97").Value = 203
$ws.Range("L97").Value = 609
$ws.Range("N97").Value = -1601

$ws.Range("H113").Value = 614.6667
$ws.Range("I113").Value = 357.44446
$ws.Range("J113").Value = 769
$ws.Range("K113").Value = 1072.33338
$ws.Range("L113").Value = 2307
$ws.Range("M113").Value = 1097.66662
$ws.Range("N113").Value = -6647


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 38463944
$ws.Range("I80").Value = 71430860
$ws.Range("J80").Value = 2534.1667
$ws.Range("K80").Value = 71430860
$ws.Range("L80").Value = 2534.1667
$ws.Range("M80").Value = -71429862
$ws.Range("N80").Value = -4530.1667

$ws.Range("H83").Value = 38463944
$ws.Range("I83").Value = 71430860
$ws.Range("J83").Value = 2534.1667
$ws.Range("K83").Value = 357154300
$ws.Range("L83").Value = 12670.8335
$ws.Range("M83").Value = -357149308
$ws.Range("N83").Value = -22654.8335

$ws.Range("H132").Value = 4736.7812
$ws.Range("I132").Value = 3822.182
$ws.Range("J132").Value = 6748.9
$ws.Range("K132").Value = 11466.546
$ws.Range("L132").Value = 20246.7
$ws.Range("M132").Value = -8936.545999999998
$ws.Range("N132").Value = -25306.7

$ws.Range("H141").Value = 76998.8
$ws.Range("J141").Value = 76666.44500000001
$ws.Range("L141").Value = 76666.44500000001
$ws.Range("N141").Value = -87026.44500000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 340.7143
$ws.Range("I9").Value = 358.33334
$ws.Range("J9").Value = 235
$ws.Range("K9").Value = 358.33334
$ws.Range("L9").Value = 235
$ws.Range("M9").Value = -134.33334
$ws.Range("N9").Value = -683

$ws.Range("H16").Value = 2161.0908
$ws.Range("I16").Value = 2161.0908
$ws.Range("K16").Value = 2161.0908
$ws.Range("M16").Value = -1991.0908

$ws.Range("H40").Value = 5558230.5
$ws.Range("I40").Value = 3106.25
$ws.Range("J40").Value = 27778728
$ws.Range("K40").Value = 3106.25
$ws.Range("L40").Value = 27778728
$ws.Range("M40").Value = -2970.25
$ws.Range("N40").Value = -27779000

$ws.Range("H132").Value = 2350
$ws.Range("I132").Value = 1363.2727
$ws.Range("J132").Value = 3706.75
$ws.Range("K132").Value = 4089.8181
$ws.Range("L132").Value = 11120.25
$ws.Range("M132").Value = -1559.8181
$ws.Range("N132").Value = -16180.25

$ws.Range("H140").Value = 73049.625
$ws.Range("J140").Value = 72913.86
$ws.Range("L140").Value = 72913.86
$ws.Range("N140").Value = -83273.86

$ws.Range("H141").Value = 110094.664
$ws.Range("J141").Value = 109913.6
$ws.Range("L141").Value = 109913.6
$ws.Range("N141").Value = -120273.6


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 184082.12
$ws.Range("J46").Value = 184082.12
$ws.Range("L46").Value = 184082.12
$ws.Range("N46").Value = -184544.12

$ws.Range("H56").Value = 23333
$ws.Range("I56").Value = 21999
$ws.Range("J56").Value = 24000
$ws.Range("K56").Value = 21999
$ws.Range("L56").Value = 24000
$ws.Range("M56").Value = -21285
$ws.Range("N56").Value = -25428

$ws.Range("H81").Value = 8792.857
$ws.Range("J81").Value = 11570
$ws.Range("L81").Value = 23140
$ws.Range("N81").Value = -25262

$ws.Range("H84").Value = 8792.857
$ws.Range("J84").Value = 11570
$ws.Range("L84").Value = 115700
$ws.Range("N84").Value = -126308

$ws.Range("H98").Value = 9000
$ws.Range("J98").Value = 9000
$ws.Range("L98").Value = 9000
$ws.Range("N98").Value = -14990

$ws.Range("H102").Value = 63000
$ws.Range("J102").Value = 63000
$ws.Range("L102").Value = 63000
$ws.Range("N102").Value = -69490

$ws.Range("H106").Value = 49000
$ws.Range("J106").Value = 49000
$ws.Range("L106").Value = 49000
$ws.Range("N106").Value = -51524

$ws.Range("H107").Value = 2429.125
$ws.Range("I107").Value = 873.6667
$ws.Range("J107").Value = 3362.4
$ws.Range("K107").Value = 2621.0001
$ws.Range("L107").Value = 10087.2
$ws.Range("M107").Value = -701.0001000000002
$ws.Range("N107").Value = -13927.2

$ws.Range("H132").Value = 1012618.25
$ws.Range("J132").Value = 4833277
$ws.Range("L132").Value = 14499831
$ws.Range("N132").Value = -14504891

$ws.Range("H134").Value = 184082.12
$ws.Range("J134").Value = 184082.12
$ws.Range("L134").Value = 552246.36
$ws.Range("N134").Value = -557316.36

$ws.Range("H137").Value = 147849.42
$ws.Range("J137").Value = 147849.42
$ws.Range("L137").Value = 147849.42
$ws.Range("N137").Value = -158049.42

$ws.Range("H140").Value = 97748.75
$ws.Range("J140").Value = 97748.75
$ws.Range("L140").Value = 97748.75
$ws.Range("N140").Value = -108108.75

